$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.921.87'
$ws.Range('E2').Value = '  -1.84%  '
$ws.Range('D3').Value = '3.138.21'
$ws.Range('E3').Value = '  -7.55%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '566.87'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -2.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.34'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -6.24%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.610'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.78%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '3.135.75'
$ws.Range('E9').Value = '  -7.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.123'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -5.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.52'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -5.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.389'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.16%  '
$ws.Range('D13').Value = '3.686.06'
$ws.Range('E13').Value = '  -7.45%  '
$ws.Range('E14').Value = '  +1.09%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.85'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -7.49%  '
$ws.Range('D16').Value = '64.820.39'
$ws.Range('E16').Value = '  -2.15%  '
$ws.Range('E17').Value = '  -6.17%  '
$ws.Range('D18').Value = '3.148.26'
$ws.Range('E18').Value = '  -7.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.69'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -3.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.73'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -7.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '355.53'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.20'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -4.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  -5.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.493'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -6.79%  '
$ws.Range('D26').Value = '3.296.76'
$ws.Range('E26').Value = '  -7.48%  '
$ws.Range('E27').Value = '  -7.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.62'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.52%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.175'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.89'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.68%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.80'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -5.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.27'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -8.08%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.19'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -4.53%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.57'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.22'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('E38').Value = '  -7.00%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.831'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.41%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.76'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '25.88'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.00%  '
$ws.Range('D42').Value = '2.659.13'
$ws.Range('E42').Value = '  -1.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.43'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -6.75%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.04'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -3.19%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.14'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -4.78%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.42'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0650'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '23.86'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -3.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '319.10'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -3.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0270'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -4.84%  '
$ws.Range('E51').Value = '  -1.52%  '
